$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before I (shifts old I:K -> J:L)
$ws.Columns("I").Insert()

# slrtype changed from Interventional to Clinical
$ws.Range("C2").Value = "Clinical"
$ws.Range("D2").Value = "Clinical_radio_button"

# The Excel/Word report filename strings (now in column J after the shift)
$ws.Range("J3").Value = "ExcelReport-NewImportLogic_1 - Test_Automation_1-Clinical-"
$ws.Range("J4").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Clinical-"

# New column I content (header + value)
$ws.Range("I1").Value = "ExpectedSourceTemplateFile"
$ws.Range("I2").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_Source_Data_Manipulated.xlsx"

# Approximate width match for the newly inserted column (engine quantizes to nearest pixel)
$ws.Columns("I").ColumnWidth = 27.67

# Update selection to match the saved view
$ws.Range("I2").Select()
